$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuesday 2:00-3:00 and 3:00-4:00 slots (rows 15 & 16) held the "random sound
# piece" show. The commit fixes its slot / player: the show title grows a
# "(...)" suffix and the artist/player is corrected to "PWUPPY PWINCESS ".
$ws.Range("C15").Value = "PWUPPY PWINCESS "
$ws.Range("C16").Value = "PWUPPY PWINCESS "
$ws.Range("B15").Value = "DIRTY STINKY SMELLY (…)"
$ws.Range("B16").Value = "DIRTY STINKY SMELLY (…)"

# The longer title needs the cell to wrap, matching the other wrapped title
# cells in this column.
$ws.Range("B15").WrapText = $true
$ws.Range("B16").WrapText = $true

# Scroll the sheet back up to where the edit happened and leave the
# selection on the edited cell.
$ws.Range("B15").Select()
$excel.ActiveWindow.ScrollRow = 7
